$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manageBuilding")

# New header cells - set in column-major order to match shared-string insertion order
$ws.Range("C1").Value = "NewBuildingName"
$ws.Range("C2").Value = "building6"
$ws.Range("C3").Value = "building10"

$ws.Range("D1").Value = "NewFloorname"
$ws.Range("D2").Value = "floorNo12"
$ws.Range("D3").Value = "FloorNo5"

# Match formatting (yellow header fill) of existing header cells
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths (closest representable values to the target 17.109375 / 17.88671875;
# this runtime quantizes ColumnWidth to 1/6-character pixel steps)
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws.Columns.Item(4).ColumnWidth = 17

# Selection
$ws.Range("B1").Select()
